$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $pp = $doc.Paragraphs.Item($i)
        if ($pp.Range.Text.Contains($text)) {
            return $pp
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Hunk 1: the "Software de Ticket para documentação" paragraph becomes
# the new "Projeto: ..." text, followed by a blank paragraph and a new
# "Software de gestão da informação e Comunicação empresarial" paragraph.
# ---------------------------------------------------------------------
$pIntro = Find-ParagraphByText $d "Software de Ticket para documentação"
$pIntro.Range.Text = "Projeto: Software de CRM (Customer Relationship Management) ou gestão de relacionamento com o cliente. "

$pIntro.Range.InsertParagraphAfter()
$pBlank = $pIntro.Next()
$pBlank.Range.Text = ""

$pBlank.Range.InsertParagraphAfter()
$pInfo = $pBlank.Next()
$pInfo.Range.Text = "Software de gestão da informação e Comunicação empresarial"

# ---------------------------------------------------------------------
# Hunk 2: two new "Problemas" bullet paragraphs ("Sistema Lento" and
# "Não customizável") are inserted right after "Interface pouco
# amigável", before the final (blank) paragraph - which itself gains an
# extra leading tab.
# ---------------------------------------------------------------------
$pInterface = Find-ParagraphByText $d "Interface pouco amigável"

$pInterface.Range.InsertParagraphAfter()
$pSlow = $pInterface.Next()
$pSlow.Range.Text = "`tSistema Lento"

$pSlow.Range.InsertParagraphAfter()
$pCustom = $pSlow.Next()
$pCustom.Range.Text = "`tNão customizável"

$pLast = $pCustom.Next()
$pLast.Range.InsertBefore("`t")
